$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "  -1.80%  "
    3  = "  -1.75%  "
    4  = "  +0.70%  "
    5  = "  -0.24%  "
    6  = "  +0.64%  "
    7  = "  -0.98%  "
    8  = "  -1.14%  "
    9  = "  -0.92%  "
    10 = "  -3.08%  "
    11 = "  -1.97%  "
    12 = "  +3.08%  "
    13 = "  -2.42%  "
    14 = "  -2.20%  "
    15 = "  -1.40%  "
    16 = "  +0.78%  "
    17 = "  +0.87%  "
    18 = "  -2.04%  "
    19 = "  +0.69%  "
    20 = "  -2.93%  "
    21 = "  -1.10%  "
    22 = "  -1.77%  "
    23 = "  +2.35%  "
    24 = "  +0.43%  "
    25 = "  -0.79%  "
    26 = "  -1.09%  "
    27 = "  -3.06%  "
    28 = "  +0.53%  "
    29 = "  -3.80%  "
    30 = "  -11.20%  "
    32 = "  -4.05%  "
    33 = "  -3.63%  "
    34 = "  -4.50%  "
    35 = "  -4.09%  "
    37 = "  -0.45%  "
    38 = "  -1.80%  "
    39 = "  -3.66%  "
    40 = "  -2.36%  "
    41 = "  -2.83%  "
    42 = "  -6.89%  "
    43 = "  -9.07%  "
    44 = "  -4.45%  "
    45 = "  -2.45%  "
    46 = "  +0.71%  "
    47 = "  -1.23%  "
    48 = "  -2.84%  "
    49 = "  -3.61%  "
    50 = "  -1.74%  "
    51 = "  -2.77%  "
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}
